$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 834, shifting the existing row 834 (and
# everything below it) down by one. This mirrors the source diff, where a
# new "2026/02/17" data row is inserted right after the existing
# "2026/02/16" row (row 833) and the rest of the table shifts down.
$ws.Rows("834:834").Insert()

# Populate the newly inserted row 834 with the new data point.
# Column A holds dates formatted as plain text (matching every other row
# in the sheet, which stores dates as literal strings rather than Excel
# date serials) -- a leading apostrophe forces Excel to keep the
# "2026/02/17" entry as text instead of auto-converting it to a date value.
$ws.Range("A834").Value = "'2026/02/17"
$ws.Range("B834").Value = "火"
$ws.Range("C834").Value = 5
$ws.Range("D834").Value = 201
